# Updated cryptos list on Sun Jul  2 21:05:18 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# the crypto rows on the active sheet. Column D values are written with a
# leading apostrophe so Excel stores them as literal text (matching the
# original inline-string cells) instead of re-parsing them as numbers and
# silently dropping trailing zeros (e.g. "166.10" -> 166.1, "1.000" -> 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = @{
    2  = "30.582.10"
    3  = "1.918.21"
    4  = "0.9995"
    5  = "245.31"
    6  = "0.9998"
    7  = "0.4832"
    8  = "0.2895"
    9  = "0.06804"
    10 = "112.14"
    12 = "1.919.18"
    13 = "0.07573"
    14 = "5.457"
    15 = "0.6727"
    16 = "295.59"
    17 = "30.568.36"
    18 = "13.02"
    19 = "0.000007612"
    21 = "5.512"
    22 = "2.164.55"
    23 = "1.000"
    25 = "9.492"
    26 = "166.10"
    27 = "20.31"
    28 = "2.091"
    29 = "0.1064"
    31 = "4.127"
    32 = "4.056"
    33 = "0.04980"
    34 = "0.7349"
    35 = "1.140"
    36 = "2.711"
    37 = "0.02029"
    39 = "2.023"
    40 = "109.47"
    41 = "0.4432"
    42 = "0.8685"
    43 = "5.847"
    44 = "1.000"
    45 = "69.36"
    46 = "7.227"
    47 = "48.77"
    48 = "9.260"
    49 = "0.1227"
    50 = "0.2511"
    51 = "34.78"
}

$volumeUpdates = @{
    2  = "  -0.05%  "
    3  = "  -0.28%  "
    4  = "  -0.09%  "
    5  = "  -0.74%  "
    6  = "  -0.05%  "
    7  = "  +1.86%  "
    8  = "  -0.31%  "
    9  = "  -0.58%  "
    10 = "  +6.09%  "
    11 = "  +5.91%  "
    12 = "  -0.26%  "
    13 = "  -1.64%  "
    14 = "  +1.93%  "
    15 = "  +0.17%  "
    16 = "  +1.64%  "
    17 = "  -0.12%  "
    18 = "  +0.52%  "
    19 = "  -0.25%  "
    20 = "  -0.06%  "
    21 = "  -0.67%  "
    22 = "  -0.43%  "
    23 = "  +0.05%  "
    24 = "  -0.69%  "
    25 = "  -0.43%  "
    26 = "  -0.87%  "
    27 = "  -3.91%  "
    28 = "  -1.40%  "
    29 = "  -0.77%  "
    30 = "  +2.81%  "
    31 = "  -1.26%  "
    32 = "  +0.04%  "
    33 = "  -1.06%  "
    34 = "  +0.40%  "
    35 = "  -0.51%  "
    36 = "  -0.49%  "
    38 = "  -0.05%  "
    39 = "  -1.01%  "
    40 = "  -2.45%  "
    41 = "  +0.38%  "
    42 = "  -0.56%  "
    43 = "  -1.45%  "
    44 = "  -0.03%  "
    45 = "  +2.15%  "
    46 = "  -0.94%  "
    47 = "  +0.35%  "
    48 = "  -1.30%  "
    49 = "  -1.39%  "
    50 = "  +0.58%  "
    51 = "  -0.75%  "
}

foreach ($row in $priceUpdates.Keys) {
    $ws.Cells.Item($row, 4).Value = "'" + $priceUpdates[$row]
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $volumeUpdates[$row]
}

Write-Output "Updated $($priceUpdates.Count) price cells and $($volumeUpdates.Count) volume cells"
